# Hortaliza, Vega Modelo de Temuco - Apio
# Weekly data refresh: insert two new daily records near the top (rows 449/450)
# and append two more records at the bottom (rows 566/567).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at 449-450, pushing the existing rows down ---
$ws.Rows("449:450").Insert()

# New row 449
$ws.Cells.Item(449, 1).Value = 10
$ws.Cells.Item(449, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(449, 3).Value = 'La Araucanía'
$ws.Cells.Item(449, 4).Value = 45204
$ws.Cells.Item(449, 5).Value = 9
$ws.Cells.Item(449, 6).Value = 100112017
$ws.Cells.Item(449, 7).Value = 'Apio'
$ws.Cells.Item(449, 8).Value = 'Americana (o)'
$ws.Cells.Item(449, 9).Value = 'Primera'
$ws.Cells.Item(449, 10).Value = 85
$ws.Cells.Item(449, 11).Value = 10000
$ws.Cells.Item(449, 12).Value = 10000
$ws.Cells.Item(449, 13).Value = 10000
$ws.Cells.Item(449, 14).Value = '$/caja 8 unidades'
$ws.Cells.Item(449, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(449, 16).Value = 10000
$ws.Cells.Item(449, 17).Value = 1
$ws.Cells.Item(449, 18).Value = 'Hortaliza'

# New row 450
$ws.Cells.Item(450, 1).Value = 10
$ws.Cells.Item(450, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(450, 3).Value = 'La Araucanía'
$ws.Cells.Item(450, 4).Value = 45204
$ws.Cells.Item(450, 5).Value = 9
$ws.Cells.Item(450, 6).Value = 100112017
$ws.Cells.Item(450, 7).Value = 'Apio'
$ws.Cells.Item(450, 8).Value = 'Americana (o)'
$ws.Cells.Item(450, 9).Value = 'Primera'
$ws.Cells.Item(450, 10).Value = 375
$ws.Cells.Item(450, 11).Value = 8000
$ws.Cells.Item(450, 12).Value = 9000
$ws.Cells.Item(450, 13).Value = 8333
$ws.Cells.Item(450, 14).Value = '$/docena de matas'
$ws.Cells.Item(450, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(450, 16).Value = 1389
$ws.Cells.Item(450, 17).Value = 6
$ws.Cells.Item(450, 18).Value = 'Hortaliza'

# --- Append two new rows (566, 567) at the bottom ---
$ws.Cells.Item(566, 1).Value = 10
$ws.Cells.Item(566, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(566, 3).Value = 'La Araucanía'
$ws.Cells.Item(566, 4).Value = 44789
$ws.Cells.Item(566, 5).Value = 9
$ws.Cells.Item(566, 6).Value = 100112017
$ws.Cells.Item(566, 7).Value = 'Apio'
$ws.Cells.Item(566, 8).Value = 'Americana (o)'
$ws.Cells.Item(566, 9).Value = 'Primera'
$ws.Cells.Item(566, 10).Value = 115
$ws.Cells.Item(566, 11).Value = 10000
$ws.Cells.Item(566, 12).Value = 11000
$ws.Cells.Item(566, 13).Value = 10435
$ws.Cells.Item(566, 14).Value = '$/docena de matas'
$ws.Cells.Item(566, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(566, 16).Value = 1739
$ws.Cells.Item(566, 17).Value = 6
$ws.Cells.Item(566, 18).Value = 'Hortaliza'

$ws.Cells.Item(567, 1).Value = 10
$ws.Cells.Item(567, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(567, 3).Value = 'La Araucanía'
$ws.Cells.Item(567, 4).Value = 44552
$ws.Cells.Item(567, 5).Value = 9
$ws.Cells.Item(567, 6).Value = 100112017
$ws.Cells.Item(567, 7).Value = 'Apio'
$ws.Cells.Item(567, 8).Value = 'Americana (o)'
$ws.Cells.Item(567, 9).Value = 'Primera'
$ws.Cells.Item(567, 10).Value = 125
$ws.Cells.Item(567, 11).Value = 8000
$ws.Cells.Item(567, 12).Value = 8000
$ws.Cells.Item(567, 13).Value = 8000
$ws.Cells.Item(567, 14).Value = '$/docena de matas'
$ws.Cells.Item(567, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(567, 16).Value = 1333
$ws.Cells.Item(567, 17).Value = 6
$ws.Cells.Item(567, 18).Value = 'Hortaliza'
